{"js": "// Rearranged abstract text:\n//  - Paragraph 1: \"contains\" -> \"contained\"; remove the \"As of 2021 ... Contact\n//    info@environmentaldatainitiative.org for more details.\" sentence (which\n//    spanned 3 runs) from the end of the first sentence group.\n//  - Paragraph 3 (\"Database history: ...\"): drop the trailing tab characters\n//    and append the \"As of 2021 ... for more details.\" text (moved from\n//    paragraph 1) to the end of the paragraph.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nlet firstPara = null;\nlet historyPara = null;\n\nfor (const p of paragraphs.items) {\n  if (firstPara === null && p.text.indexOf(\"This dataset is an archive of the ClimHydroDB database\") !== -1) {\n    firstPara = p;\n  }\n  if (historyPara === null && p.text.indexOf(\"Database history:\") !== -1) {\n    historyPara = p;\n  }\n}\n\nif (!firstPara || !historyPara) {\n  throw new Error(\"Could not locate the expected paragraphs in the document.\");\n}\n\nconst newFirstParaText =\n  \"This dataset is an archive of the ClimHydroDB database, which was actively used from early 2001 to mid 2020. \" +\n  \"The database contained contributions from 62 contributors (primarily from the LTER Network and US Forest Service) and 672 research sites. \" +\n  \"Data records total approximately 16 million (raw) or 1.6 million (aggregated) for 22 meteorologic or hydrologic variables. \" +\n  \"This archive contains the 23 core tables of the ClimHydroDB database as text tables of comma separated values, plus the database entity relationship diagram (ERD), User Guide, database table descriptions (DDL, SQL script), and a zip file of related documents and presentations.\";\n\nconst newHistoryParaText =\n  \"Database history: To facilitate intersite research within the LTER network, site data managers developed a system to provide climatic summaries dynamically, called ClimDB. \" +\n  \"Later funding from the U. S. Forest Service allowed the original database to be expanded to include hydrologic variables, and the combined database was renamed ClimHydroDB in 2003. \" +\n  \"The database also harvested real-time streamflow data from USGS gauging stations, using code developed by the Georgia Coastal Ecosystem LTER.  \" +\n  \"As of 2021, the ClimHydroDB content is available as data packages from individual contributing sites, each containing identically formatted text tables in the ODM 1.1 format, for integration with CUAHSI tools (https://cuahsi.org). Contact info@environmentaldatainitiative.org for more details.\";\n\n// Replacing the whole paragraph's text collapses its runs into a single run\n// that keeps the formatting (sz/szCs/rtl) of the paragraph's original first run.\nfirstPara.insertText(newFirstParaText, \"Replace\");\nhistoryPara.insertText(newHistoryParaText, \"Replace\");\n\nawait context.sync();\n", "ps1": "# Rearranged abstract text:\n#  - Paragraph 1: \"contains\" -> \"contained\"; remove the \"As of 2021 ... Contact\n#    info@environmentaldatainitiative.org for more details.\" sentence (which\n#    spanned 3 runs) from the end of the first sentence group.\n#  - Paragraph 3 (\"Database history: ...\"): drop the trailing tab characters\n#    and append the \"As of 2021 ... for more details.\" text (moved from\n#    paragraph 1) to the end of the paragraph.\n\n$d = $word.ActiveDocument\n$wdFindContinue = 1\n$wdReplaceAll = 2\n\nfunction Replace-Text($findText, $replaceText) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $findText\n    $find.Replacement.Text = $replaceText\n    $find.Execute($find.Text, $true, $false, $false, $false, $false, $true, $wdFindContinue, $false, $find.Replacement.Text, $wdReplaceAll) | Out-Null\n}\n\n# 1) \"contains\" -> \"contained\" (target the specific sentence so we don't touch\n#    any other occurrence of the word \"contains\" elsewhere in the document).\nReplace-Text \"The database contains contributions from 62 contributors\" \"The database contained contributions from 62 contributors\"\n\n# 2) Remove the \"As of 2021 ... for more details.\" sentence from paragraph 1\n#    (it moves down into the \"Database history\" paragraph below).\nReplace-Text \"variables. As of 2021, the ClimHydroDB content is available as data packages from individual contributing sites, each containing identically formatted text tables in the ODM 1.1 format, for integration with CUAHSI tools (https://cuahsi.org). Contact info@environmentaldatainitiative.org for more details. This archive\" \"variables. This archive\"\n\n# 3) Drop the trailing tabs on the \"Database history\" paragraph and append the\n#    sentence moved from paragraph 1.\n$tab = [char]9\n$findTabs = \"Georgia Coastal Ecosystem LTER.  \" + $tab + $tab + $tab + $tab\nReplace-Text $findTabs \"Georgia Coastal Ecosystem LTER.  As of 2021, the ClimHydroDB content is available as data packages from individual contributing sites, each containing identically formatted text tables in the ODM 1.1 format, for integration with CUAHSI tools (https://cuahsi.org). Contact info@environmentaldatainitiative.org for more details.\"\n"}
